# Refined metadata to be additional tab
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("data")

# --- Update the "panel_query_time" style timestamps recorded per-gene on
# --- the "data" sheet (column F) to reflect the refreshed query run time.
$ws1.Range("F2").Value = "2021-10-05 14:33:33.052689"
$ws1.Range("F3").Value = "2021-10-05 14:33:33.052697"
$ws1.Range("F4").Value = "2021-10-05 14:33:33.052701"
$ws1.Range("F5").Value = "2021-10-05 14:33:33.052704"
$ws1.Range("F6").Value = "2021-10-05 14:33:33.052707"
$ws1.Range("F7").Value = "2021-10-05 14:33:33.052710"
$ws1.Range("F8").Value = "2021-10-05 14:33:33.052713"
$ws1.Range("F9").Value = "2021-10-05 14:33:33.052716"
$ws1.Range("F10").Value = "2021-10-05 14:33:33.052719"
$ws1.Range("F11").Value = "2021-10-05 14:33:33.052722"
$ws1.Range("F12").Value = "2021-10-05 14:33:33.052724"
$ws1.Range("F13").Value = "2021-10-05 14:33:33.052727"
$ws1.Range("F14").Value = "2021-10-05 14:33:33.052730"
$ws1.Range("F15").Value = "2021-10-05 14:33:33.052732"
$ws1.Range("F16").Value = "2021-10-05 14:33:33.052735"
$ws1.Range("F17").Value = "2021-10-05 14:33:33.052738"
$ws1.Range("F18").Value = "2021-10-05 14:33:33.052741"
$ws1.Range("F19").Value = "2021-10-05 14:33:33.052744"
$ws1.Range("F20").Value = "2021-10-05 14:33:33.052746"
$ws1.Range("F21").Value = "2021-10-05 14:33:33.052749"
$ws1.Range("F22").Value = "2021-10-05 14:33:33.052752"
$ws1.Range("F23").Value = "2021-10-05 14:33:33.052755"
$ws1.Range("F24").Value = "2021-10-05 14:33:33.052758"
$ws1.Range("F25").Value = "2021-10-05 14:33:33.052760"
$ws1.Range("F26").Value = "2021-10-05 14:33:33.052764"
$ws1.Range("F27").Value = "2021-10-05 14:33:33.052767"
$ws1.Range("F28").Value = "2021-10-05 14:33:33.052769"
$ws1.Range("F29").Value = "2021-10-05 14:33:33.052772"
$ws1.Range("F30").Value = "2021-10-05 14:33:33.052775"
$ws1.Range("F31").Value = "2021-10-05 14:33:33.052777"
$ws1.Range("F32").Value = "2021-10-05 14:33:33.052780"
$ws1.Range("F33").Value = "2021-10-05 14:33:33.052783"
$ws1.Range("F34").Value = "2021-10-05 14:33:33.052786"
$ws1.Range("F35").Value = "2021-10-05 14:33:33.052788"
$ws1.Range("F36").Value = "2021-10-05 14:33:33.052791"
$ws1.Range("F37").Value = "2021-10-05 14:33:33.052794"
$ws1.Range("F38").Value = "2021-10-05 14:33:33.052796"
$ws1.Range("F39").Value = "2021-10-05 14:33:33.052799"
$ws1.Range("F40").Value = "2021-10-05 14:33:33.052802"
$ws1.Range("F41").Value = "2021-10-05 14:33:33.052805"
$ws1.Range("F42").Value = "2021-10-05 14:33:33.052808"

# --- Add the new "metadata" sheet (after "data") that records the
# --- PanelApp query parameters used to build this export.
$meta = $wb.Worksheets.Add($null, $ws1)
$meta.Name = "metadata"
$meta.Visible = -1

# Header row (bold / bordered / centered) mirrors the "data" sheet's header
# style, so copy formats from an existing header cell onto the new header.
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$ws1.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

# Data row
$meta.Range("A2").Value = 0
$ws1.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

$meta.Range("B2").Value = "Congenital anomalies of the kidney and urinary tract (CAKUT) Nonsyndromic"
$meta.Range("C2").Value = 61

# data_version ("0.89") must be stored as text, not coerced to a number -
# format the cell as Text before entry, then restore the default style so
# the cell doesn't carry a stray number-format flag.
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "0.89"
$meta.Range("D2").Style = "Normal"

$meta.Range("E2").Value = "2021-09-17T09:54:15.404193Z"
$meta.Range("F2").Value = "2021-10-05 14:33:33.048825"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/61/?format=json"

$excel.CutCopyMode = $false
